$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B7").Value = "sogou"
$ws.Range("D7").Value = 1567
$ws.Range("D8").Value = 6508
$ws.Range("D9").Value = 5873
$ws.Range("D10").Value = 1854
